$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Qty executed upto date" (column C) numeric updates ---
$ws.Range("C8").Value = 51
$ws.Range("C9").Value = 47
$ws.Range("C10").Value = 17
$ws.Range("C11").Value = 2
$ws.Range("C12").Value = 13
$ws.Range("C13").Value = 53
$ws.Range("C14").Value = 12
$ws.Range("C15").Value = 80
$ws.Range("C16").Value = 54
$ws.Range("C17").Value = 12

# --- "Upto date Amount" (column G) text-formatted updates ---
# These cells store formatted currency amounts as text (e.g. "12032.00"),
# so force Text formatting before assigning the string value to keep
# Excel from re-interpreting it as a number.
$gUpdates = @{
    "G9"  = "12032.00"
    "G10" = "8024.00"
    "G11" = "1324.00"
    "G13" = "7208.00"
    "G14" = "276.00"
    "G19" = "28864.00"
    "G21" = "28864.00"
}
foreach ($addr in $gUpdates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $gUpdates[$addr]
}

# --- "Amount Since prev bill" (column H) grand-total text updates ---
$hUpdates = @{
    "H19" = "28864.00"
    "H21" = "28864.00"
}
foreach ($addr in $hUpdates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $hUpdates[$addr]
}
